# Auto-generated script: applies cached-value refresh to Marilith_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5919.5
$ws.Range("J40").Value = 6705.294
$ws.Range("L40").Value = 6705.294
$ws.Range("N40").Value = -7055.294
$ws.Range("H112").Value = 2330.7693
$ws.Range("J112").Value = 2720
$ws.Range("L112").Value = 8160
$ws.Range("N112").Value = -10376
$ws.Range("H115").Value = 82
$ws.Range("I115").Value = 82
$ws.Range("K115").Value = 246
$ws.Range("M115").Value = 1321
$ws.Range("H132").Value = 2711.7693
$ws.Range("I132").Value = 3059.3635
$ws.Range("K132").Value = 9178.0905
$ws.Range("M132").Value = -6648.0905
$ws.Range("H138").Value = 2074.8462
$ws.Range("I138").Value = 994.6
$ws.Range("K138").Value = 2983.8
$ws.Range("M138").Value = 2156.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 810.36365
$ws.Range("I45").Value = 819.9
$ws.Range("J45").Value = 715
$ws.Range("K45").Value = 819.9
$ws.Range("L45").Value = 715
$ws.Range("M45").Value = -442.9
$ws.Range("N45").Value = -1469
$ws.Range("H61").Value = 3860.1428
$ws.Range("I61").Value = 2836.8333
$ws.Range("K61").Value = 2836.8333
$ws.Range("M61").Value = -2624.8333
$ws.Range("H106").Value = 22000
$ws.Range("J106").Value = 22000
$ws.Range("L106").Value = 22000
$ws.Range("N106").Value = -24524
$ws.Range("H122").Value = 2819.6
$ws.Range("I122").Value = 2727.3572
$ws.Range("K122").Value = 8182.071599999999
$ws.Range("M122").Value = -5732.071599999999
$ws.Range("H136").Value = 3860.1428
$ws.Range("I136").Value = 2836.8333
$ws.Range("K136").Value = 8510.499899999999
$ws.Range("M136").Value = -5960.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 820.2857
$ws.Range("I20").Value = 638.125
$ws.Range("J20").Value = 1063.1666
$ws.Range("K20").Value = 638.125
$ws.Range("L20").Value = 1063.1666
$ws.Range("M20").Value = -391.125
$ws.Range("N20").Value = -1557.1666
$ws.Range("H29").Value = 733
$ws.Range("I29").Value = 733
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 733
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -444
$ws.Range("H86").Value = 3488.5454
$ws.Range("I86").Value = 4479.8
$ws.Range("J86").Value = 2662.5
$ws.Range("K86").Value = 4479.8
$ws.Range("L86").Value = 2662.5
$ws.Range("M86").Value = -3356.8
$ws.Range("N86").Value = -4908.5
$ws.Range("H89").Value = 3488.5454
$ws.Range("I89").Value = 4479.8
$ws.Range("J89").Value = 2662.5
$ws.Range("K89").Value = 22399
$ws.Range("L89").Value = 13312.5
$ws.Range("M89").Value = -16783
$ws.Range("N89").Value = -24544.5
$ws.Range("H94").Value = 870
$ws.Range("I94").Value = 865
$ws.Range("K94").Value = 865
$ws.Range("M94").Value = -414
$ws.Range("H95").Value = 39000
$ws.Range("J95").Value = 39000
$ws.Range("L95").Value = 39000
$ws.Range("N95").Value = -44492
$ws.Range("H100").Value = 19721.334
$ws.Range("J100").Value = 19721.334
$ws.Range("L100").Value = 19721.334
$ws.Range("N100").Value = -21885.334
$ws.Range("H105").Value = 3018.8
$ws.Range("I105").Value = 3098.6667
$ws.Range("K105").Value = 3098.6667
$ws.Range("M105").Value = -1351.6667
$ws.Range("H134").Value = 5247.8076
$ws.Range("I134").Value = 5192.875
$ws.Range("K134").Value = 15578.625
$ws.Range("M134").Value = -13043.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 456654.47
$ws.Range("I35").Value = 502289.9
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 502289.9
$ws.Range("L35").Value = 300
$ws.Range("M35").Value = -501995.9
$ws.Range("N35").Value = -888
$ws.Range("H43").Value = 12828.5
$ws.Range("J43").Value = 12828.5
$ws.Range("L43").Value = 12828.5
$ws.Range("N43").Value = -13196.5
$ws.Range("H88").Value = 32666.5
$ws.Range("J88").Value = 32666.5
$ws.Range("L88").Value = 32666.5
$ws.Range("N88").Value = -33478.5
$ws.Range("H91").Value = 32666.5
$ws.Range("J91").Value = 32666.5
$ws.Range("L91").Value = 32666.5
$ws.Range("N91").Value = -35474.5
$ws.Range("H101").Value = 12828.5
$ws.Range("J101").Value = 12828.5
$ws.Range("L101").Value = 12828.5
$ws.Range("N101").Value = -19318.5
$ws.Range("H132").Value = 4113.2856
$ws.Range("I132").Value = 4113.2856
$ws.Range("K132").Value = 12339.8568
$ws.Range("M132").Value = -9809.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2204.4285
$ws.Range("J4").Value = 2079.5
$ws.Range("L4").Value = 6238.5
$ws.Range("N4").Value = -6462.5
$ws.Range("H52").Value = 1749.5
$ws.Range("J52").Value = 1749.5
$ws.Range("L52").Value = 5248.5
$ws.Range("N52").Value = -5780.5
$ws.Range("H55").Value = 365
$ws.Range("J55").Value = 325
$ws.Range("L55").Value = 975
$ws.Range("N55").Value = -1329
$ws.Range("H131").Value = 1719.6
$ws.Range("I131").Value = 866.3333
$ws.Range("K131").Value = 2598.9999
$ws.Range("M131").Value = 2441.0001
$ws.Range("H139").Value = 1666
$ws.Range("I139").Value = 999
$ws.Range("K139").Value = 2997
$ws.Range("M139").Value = 2143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 14000
$ws.Range("J33").Value = 14000
$ws.Range("L33").Value = 14000
$ws.Range("N33").Value = -14504
$ws.Range("H80").Value = 2452.2222
$ws.Range("J80").Value = 2503
$ws.Range("L80").Value = 2503
$ws.Range("N80").Value = -4499
$ws.Range("H83").Value = 2452.2222
$ws.Range("J83").Value = 2503
$ws.Range("L83").Value = 12515
$ws.Range("N83").Value = -22499
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
$ws.Range("H102").Value = 953.0833
$ws.Range("I102").Value = 793.1111
$ws.Range("J102").Value = 1433
$ws.Range("K102").Value = 793.1111
$ws.Range("L102").Value = 1433
$ws.Range("M102").Value = 828.8889
$ws.Range("N102").Value = -4677
$ws.Range("H122").Value = 1095.2858
$ws.Range("I122").Value = 1139
$ws.Range("K122").Value = 3417
$ws.Range("M122").Value = -967

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12853.5
$ws.Range("I7").Value = 12662.417
$ws.Range("K7").Value = 12662.417
$ws.Range("M7").Value = -12550.417
$ws.Range("H76").Value = 34572.8
$ws.Range("J76").Value = 34572.8
$ws.Range("L76").Value = 34572.8
$ws.Range("N76").Value = -35248.8
$ws.Range("H79").Value = 34572.8
$ws.Range("J79").Value = 34572.8
$ws.Range("L79").Value = 34572.8
$ws.Range("N79").Value = -36912.8
$ws.Range("H101").Value = 21454
$ws.Range("J101").Value = 21454
$ws.Range("L101").Value = 21454
$ws.Range("N101").Value = -27944
$ws.Range("H119").Value = 252000
$ws.Range("J119").Value = 252000
$ws.Range("L119").Value = 252000
$ws.Range("N119").Value = -261676
$ws.Range("H126").Value = 12853.5
$ws.Range("I126").Value = 12662.417
$ws.Range("K126").Value = 37987.251
$ws.Range("M126").Value = -35517.251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 22999.5
$ws.Range("J69").Value = 22999.5
$ws.Range("L69").Value = 22999.5
$ws.Range("N69").Value = -24497.5
$ws.Range("H72").Value = 22999.5
$ws.Range("J72").Value = 22999.5
$ws.Range("L72").Value = 68998.5
$ws.Range("N72").Value = -76486.5
$ws.Range("H81").Value = 625
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 625
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
$ws.Range("H132").Value = 2278.75
$ws.Range("I132").Value = 2094.111
$ws.Range("K132").Value = 6282.333
$ws.Range("M132").Value = -3752.333
